$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.046.42"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.831.09"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'241.55"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "'0.6268"
$ws.Range("E6").Value = "  -5.06%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.07609"
$ws.Range("E8").Value = "  +2.54%  "
$ws.Range("D9").Value = "'0.2917"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "'22.80"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "'0.07640"
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("D12").Value = "1.833.17"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").Value = "'0.6652"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "'82.36"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").Value = "'0.000009497"
$ws.Range("E16").Value = "  +10.91%  "
$ws.Range("D17").Value = "'5.989"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("D18").Value = "28.965.37"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").Value = "'225.06"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'7.219"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'161.11"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").Value = "'8.417"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").Value = "'0.1363"
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("D27").Value = "'17.84"
$ws.Range("D28").Value = "'1.492"
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'4.036"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'4.055"
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("D31").Value = "'1.195"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("D32").Value = "'0.05199"
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("D33").Value = "'1.850"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "'0.7288"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("E36").Value = "  -2.15%  "
$ws.Range("D37").Value = "1.273.94"
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "'0.01786"
$ws.Range("D40").Value = "'6.509"
$ws.Range("E40").Value = "  +7.84%  "
$ws.Range("D41").Value = "'0.8915"
$ws.Range("E41").Value = "  -3.03%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "'101.40"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").Value = "1.975.08"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("D46").Value = "'63.73"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "'0.3980"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").Value = "'0.07328"
$ws.Range("E49").Value = "  -12.05%  "
$ws.Range("D50").Value = "'8.840"
$ws.Range("E50").Value = "  +1.35%  "
$ws.Range("D51").Value = "'0.05749"
$ws.Range("E51").Value = "  -1.58%  "

# Reset style on cells where we had to force text via leading apostrophe,
# so we do not leave a stray quote-prefix style on an otherwise default-styled cell.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
